# Bulk upload template fix:
#  - rename "ResourceCreator.*" header labels to distinguish the
#    Institution-creator bean from the Person-creator bean
#    (ResourceCreatorInstitution.* / ResourceCreatorPerson.*)
#  - move the active selection from A4 to E2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E1").Value = "ResourceCreatorInstitution.Institution.name"
$ws.Range("F1").Value = "ResourceCreatorInstitution.role"
$ws.Range("G1").Value = "ResourceCreatorPerson.Person.lastName"
$ws.Range("H1").Value = "ResourceCreatorPerson.Person.firstName"
$ws.Range("I1").Value = "ResourceCreatorPerson.Person.email"
$ws.Range("J1").Value = "ResourceCreatorPerson.Person.Institution.name"
$ws.Range("K1").Value = "ResourceCreatorPerson.role"

$ws.Select()
$ws.Range("E2").Select()
